$d = $word.ActiveDocument

# 1. Re-style the document title from the "Title" style to "Heading 1".
#    Applying the new paragraph style also clears the paragraph's
#    direct formatting overrides (left indent / centered alignment)
#    that had been layered on top of the old "Title" style.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.ParagraphFormat.Style = "Heading 1"

# 2. The sentence "...you must follow to connect to OsGrid " had been
#    split across three runs so a grammar-check proofing mark
#    (gramStart/gramEnd) could wrap the word "OsGrid". Re-running
#    Find & Replace over that exact text collapses it back into a
#    single contiguous run and drops the now-unneeded proofErr marks.
$d.Content.Find.Execute(
    "to connect to OsGrid ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "to connect to OsGrid ", 2
) | Out-Null
